$ErrorActionPreference = "Stop"

$d = $word.ActiveDocument

# Sanity check: locate the target paragraph (the "第四段" Q&A item) via Find,
# then confirm it is Paragraphs.Item(23) before mutating anything.
$checkRange = $d.Content
$found = $checkRange.Find.Execute("【問】麻烦你帮我白话翻译一下『釋四伏難』四段文", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not locate anchor paragraph for the edit."
}

$targetPara = $d.Paragraphs.Item(23)
if ($targetPara.Range.Start -ne $checkRange.Start) {
    throw "Paragraph index assumption (23) does not match the located anchor text."
}

# 1) Add w:hint="eastAsia" to the target paragraph's own paragraph-mark run
#    properties (<w:pPr><w:rPr><w:rFonts .../>). We do this by replacing the
#    whole paragraph's Range contents with an identical OOXML paragraph that
#    carries the extra w:hint attribute on that one <w:rFonts/> element.
$xmlPara23 = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="a3"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:adjustRightInd w:val="0"/><w:spacing w:beforeLines="50"/><w:ind w:left="993" w:hanging="993"/><w:rPr><w:rFonts w:ascii="华文中宋" w:eastAsia="华文中宋" w:hAnsi="华文中宋" w:hint="eastAsia"/><w:snapToGrid w:val="0"/><w:kern w:val="0"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="华文中宋" w:eastAsia="华文中宋" w:hAnsi="华文中宋" w:hint="eastAsia"/><w:snapToGrid w:val="0"/><w:kern w:val="0"/></w:rPr><w:t>【問】麻烦你帮我白话翻译一下『釋四伏難』四段文？【答】第一段，弥勒请文殊</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="华文中宋" w:eastAsia="华文中宋" w:hAnsi="华文中宋" w:hint="eastAsia"/><w:snapToGrid w:val="0"/><w:kern w:val="0"/></w:rPr><w:lastRenderedPageBreak/><w:t>回答大众疑惑。文殊内心反驳说大众并没有疑惑，如果大众有疑惑应该发问才是。大众既没有疑惑，我又何须答疑呢？弥勒就说，大众都在看着您与我。看着我，是希望我能够发问；看着您，是希望您能够回答。第二段，文殊内心说，既然大众都有疑惑，那就不容易回答。不如等佛出定后亲自回答，大众自然明白。弥勒催促说，既然大众疑惑、忧心，应该及时回答。还不知道佛什么时候会出定呢？第三段，文殊内心说，我与您两个都是佛的学生，想要猜测老师佛的意思，应该共同细心思维，让我一个人回答，道理上也说不过去。弥勒就说，我确实也在细心思维，进退维谷，是佛将要宣讲《妙法莲华经》？是佛将要给大众弟子授未来成佛之记？第四段，文殊内心说，如您所说，那大众的疑惑就已经解除了，何须我再回答什么呢？弥勒说，岂能以我这种没有十分把握的猜测而判说这样的大事因缘？文殊内心不再反驳弥勒，弥勒就乘机说，文殊您应该明白，大众都在看着您、望您答疑呢，接下来佛该说何经法呢？</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$targetPara.Range.InsertXML($xmlPara23)

# 2) Insert a brand-new Q&A paragraph right after it with the three runs
#    (plain lead-in, bold+underlined classical quotation, plain answer).
$targetPara = $d.Paragraphs.Item(23)
$insertionPoint = $targetPara.Range
$insertionPoint.Collapse(0)
$insertionPoint.InsertParagraphAfter()
$newPara = $d.Paragraphs.Item(24)
$xmlNewPara = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="a3"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:adjustRightInd w:val="0"/><w:spacing w:beforeLines="50"/><w:ind w:left="993" w:hanging="993"/><w:rPr><w:rFonts w:ascii="华文中宋" w:eastAsia="华文中宋" w:hAnsi="华文中宋"/><w:snapToGrid w:val="0"/><w:kern w:val="0"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="华文中宋" w:eastAsia="华文中宋" w:hAnsi="华文中宋" w:hint="eastAsia"/><w:snapToGrid w:val="0"/><w:kern w:val="0"/></w:rPr><w:t>【問】请翻译一下『</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="华文中宋" w:eastAsia="华文中宋" w:hAnsi="华文中宋" w:hint="eastAsia"/><w:b/><w:snapToGrid w:val="0"/><w:kern w:val="0"/><w:u w:val="single"/></w:rPr><w:t>夫以下測上，止可罔像卜度，惟昔儔今 ，不可頓決，所以初從髣髴；次引略見；略見未周，更引廣見；以多證一，爾乃分判。</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="华文中宋" w:eastAsia="华文中宋" w:hAnsi="华文中宋" w:hint="eastAsia"/><w:snapToGrid w:val="0"/><w:kern w:val="0"/></w:rPr><w:t>』一段。【答】以文殊等觉位，来猜测佛妙觉位，以下位测上位，只能够依稀仿佛，猜个大概而已。回忆往昔旧事，而来比较目前所发生之事，不能一下子决定是否如此。所以文殊从一开始的大概猜测，到随后的回忆过去略略曾经见到过相似的情景。略略曾经见过的情景还不够详细，文殊又再回忆了过去详细的相似的情景。文殊以智慧推断，并以往昔时的详略不同的多种见闻，来作证明，确定现在释迦佛确实如同往昔时的日月灯明佛一样，即将宣讲《妙法莲华经》。</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$newPara.Range.InsertXML($xmlNewPara)

Write-Output "DONE"
Write-Output $d.Paragraphs.Count
